# Apply the "Update subtype value define" change to doc/diagram.xlsx
#
# Summary of the change (see commit message / diff):
#  - New table of "flow of adding new packet type" steps written into
#    worksheet cells F167:G174 (values + shared strings).
#  - Three existing diagram boxes get renamed/re-worded:
#      my_logging.py    -> ctrl   ("Control packet")
#      my_matplotlib.py -> data   ("Data packet")
#      my_time.py       -> mgt    ("Managetment packet")
#  - A new small diagram group is added summarizing
#    "ctrl/data/mgt" with "packet_type.py" / "layer_name.py" boxes.
#  - The sheet view is refreshed (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------
# 1. New "flow of adding new packet type" mini table (F167:G174)
# ---------------------------------------------------------------------

$ws.Range("G168").Value = "In config.ini, define enable flag"
$ws.Range("F167").Value = "Flow of adding new packet type"
$ws.Range("G172").Value = "Create /my_sniff/subtype/packet_type folder, and in this folder, create packet_type.py file which includes all fields and value define"
$ws.Range("G169").Value = "In  /my_config/config_basic.py,  add func to get enable flag"
$ws.Range("G170").Value = "In /packet_config/ folder, create config_xxx.ini (xxx is the packet type name), and update related fields. In /my_config/ folder, add config_xxx.py to get values from ini file"
$ws.Range("G171").Value = "in my_flow.py, import config_xxx.py file"
$ws.Range("G173").Value = "in my_flow.py, import packet_type.py file"
$ws.Range("G174").Value = "update my_flow.py to add this packet type check section"

$ws.Range("F168").Value = 0
$ws.Range("F169").Value = 0.1
$ws.Range("F170").Value = 1
$ws.Range("F171").Value = 1.1
$ws.Range("F172").Value = 2
$ws.Range("F173").Value = 2.1
$ws.Range("F174").Value = 3

# ---------------------------------------------------------------------
# 2. Re-word the three existing "ctrl / data / mgt" related boxes
# ---------------------------------------------------------------------

# my_logging.py -> ctrl
$ctrlShape = $ws.Shapes.Item($ws.Shapes.Count - 2)
$ctrlShape.TextFrame2.TextRange.Text = "ctrl" + $nl + "---------------------------------------------------" + $nl + "Control packet" + $nl + "---------------------------------------------------" + $nl + $nl

# my_matplotlib.py -> data
$dataShape = $ws.Shapes.Item($ws.Shapes.Count - 1)
$dataShape.TextFrame2.TextRange.Text = "data" + $nl + "---------------------------------------------------" + $nl + "Data packet" + $nl + "---------------------------------------------------" + $nl + $nl

# my_time.py -> mgt
$mgtShape = $ws.Shapes.Item($ws.Shapes.Count)
$mgtShape.TextFrame2.TextRange.Text = "mgt" + $nl + "---------------------------------------------------" + $nl + "Managetment packet" + $nl + "---------------------------------------------------" + $nl

# ---------------------------------------------------------------------
# 3. New summary group: "ctrl/data/mgt" with packet_type.py / layer_name.py
# ---------------------------------------------------------------------

$fromCell = $ws.Cells.Item(151, 28)
$toCell = $ws.Cells.Item(165, 38)
$groupLeft = $fromCell.Left + (92777 / 12700)
$groupTop = $fromCell.Top + (105146 / 12700)
$groupRight = $toCell.Left + (415407 / 12700)
$groupBottom = $toCell.Top + (181346 / 12700)
$groupWidth = $groupRight - $groupLeft
$groupHeight = $groupBottom - $groupTop

$bigRect = $ws.Shapes.AddShape(5, $groupLeft, $groupTop, $groupWidth, $groupHeight)
$bigRect.Name = "Group 86"
$bigRect.TextFrame2.TextRange.Text = "ctrl/data/mgt"
$bigRect.Fill.ForeColor.RGB = 0x00A03070

$packetTypeLeft = $groupLeft + 0.0448178 * $groupWidth
$packetTypeTop = $groupTop + 0.1552286 * $groupHeight
$packetTypeWidth = 0.4285714 * $groupWidth
$packetTypeHeight = 0.6666667 * $groupHeight
$packetTypeShape = $ws.Shapes.AddShape(5, $packetTypeLeft, $packetTypeTop, $packetTypeWidth, $packetTypeHeight)
$packetTypeShape.Name = "Rounded Rectangle 5"
$packetTypeShape.TextFrame2.TextRange.Text = "packet_type.py" + $nl + "-------------------------------------------------------" + $nl + $nl + "-------------------------------------------------------"

$layerNameLeft = $groupLeft + 0.5182071 * $groupWidth
$layerNameTop = $packetTypeTop
$layerNameWidth = $packetTypeWidth
$layerNameHeight = $packetTypeHeight
$layerNameShape = $ws.Shapes.AddShape(5, $layerNameLeft, $layerNameTop, $layerNameWidth, $layerNameHeight)
$layerNameShape.Name = "Rounded Rectangle 4"
$layerNameShape.TextFrame2.TextRange.Text = "layer_name.py" + $nl + "-------------------------------------------------------" + $nl + $nl + "-------------------------------------------------------"

# ---------------------------------------------------------------------
# 4. Refresh sheet view (zoom + active selection)
# ---------------------------------------------------------------------

$excel.ActiveWindow.Zoom = 55
$ws.Range("AG23").Select()
